# Workbook: public/mau.xlsx — "xong moi chuc nang co ban"
#
# The template's sample data is trimmed down and two of the remaining
# "Đơn Vị" (unit) entries are corrected:
#   - C8 ("Phùng Thị Tuyết Lan") unit changes from "Hồ Chí Minh" to "HCM"
#   - C10 ("Chim Sẻ Đi Nắng") unit changes from "Đế Chế" to "Hà Nội"
# Rows 11-20 (the remaining sample participants 8-16 plus a spacer row)
# are removed entirely, shrinking the sheet from 30 rows to 20 rows.
# Finally the view is left scrolled near the top with C8 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the two "Đơn Vị" values that changed ---
$ws.Range("C8").Value = "HCM"
$ws.Range("C10").Value = "Hà Nội"

# --- Remove the now-unneeded sample rows (entire row delete, shifts
#     the trailing blank rows up so the sheet ends at row 20) ---
$ws.Rows("11:20").Delete()

# --- Leave the window scrolled/selected the way the author left it ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select()
